$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = '@'
$cell.Value = '309.92'
$cell.Style = 'Normal'

$cell = $ws.Range("E2")
$cell.NumberFormat = '@'
$cell.Value = '0.28%'
$cell.Style = 'Normal'

$cell = $ws.Range("D3")
$cell.NumberFormat = '@'
$cell.Value = '41.04'
$cell.Style = 'Normal'

$cell = $ws.Range("E3")
$cell.NumberFormat = '@'
$cell.Value = '-0.55%'
$cell.Style = 'Normal'

$cell = $ws.Range("D4")
$cell.NumberFormat = '@'
$cell.Value = '5.217'
$cell.Style = 'Normal'

$cell = $ws.Range("E4")
$cell.NumberFormat = '@'
$cell.Value = '1.62%'
$cell.Style = 'Normal'

$cell = $ws.Range("D5")
$cell.NumberFormat = '@'
$cell.Value = '0.07680'
$cell.Style = 'Normal'

$ws.Range("B6").Value = 'GateToken'

$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$cell = $ws.Range("D6")
$cell.NumberFormat = '@'
$cell.Value = '4.291'
$cell.Style = 'Normal'

$cell = $ws.Range("E6")
$cell.NumberFormat = '@'
$cell.Value = '0.57%'
$cell.Style = 'Normal'

$ws.Range("B7").Value = 'FTXToken'

$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

$cell = $ws.Range("D7")
$cell.NumberFormat = '@'
$cell.Value = '1.705'
$cell.Style = 'Normal'

$cell = $ws.Range("E7")
$cell.NumberFormat = '@'
$cell.Value = '5.38%'
$cell.Style = 'Normal'

$ws.Range("B8").Value = 'MXToken'

$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$cell = $ws.Range("D8")
$cell.NumberFormat = '@'
$cell.Value = '0.9385'
$cell.Style = 'Normal'

$cell = $ws.Range("E8")
$cell.NumberFormat = '@'
$cell.Value = '3.33%'
$cell.Style = 'Normal'

$ws.Range("B9").Value = 'BTSEToken'

$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$cell = $ws.Range("D9")
$cell.NumberFormat = '@'
$cell.Value = '2.425'
$cell.Style = 'Normal'

$cell = $ws.Range("E9")
$cell.NumberFormat = '@'
$cell.Value = '-1.94%'
$cell.Style = 'Normal'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$cell = $ws.Range("D10")
$cell.NumberFormat = '@'
$cell.Value = '0.1281'
$cell.Style = 'Normal'

$cell = $ws.Range("E10")
$cell.NumberFormat = '@'
$cell.Value = '10.88%'
$cell.Style = 'Normal'

$ws.Range("B11").Value = 'WazirX'

$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$cell = $ws.Range("D11")
$cell.NumberFormat = '@'
$cell.Value = '0.1834'
$cell.Style = 'Normal'

$cell = $ws.Range("E11")
$cell.NumberFormat = '@'
$cell.Value = '1.54%'
$cell.Style = 'Normal'

$ws.Range("B12").Value = 'MandalaExchangeToken'

$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$cell = $ws.Range("D12")
$cell.NumberFormat = '@'
$cell.Value = '0.09146'
$cell.Style = 'Normal'

$cell = $ws.Range("E12")
$cell.NumberFormat = '@'
$cell.Value = '-0.15%'
$cell.Style = 'Normal'

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$cell = $ws.Range("D13")
$cell.NumberFormat = '@'
$cell.Value = '0.04229'
$cell.Style = 'Normal'

$cell = $ws.Range("E13")
$cell.NumberFormat = '@'
$cell.Value = '-0.47%'
$cell.Style = 'Normal'

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$cell = $ws.Range("D14")
$cell.NumberFormat = '@'
$cell.Value = '0.1052'
$cell.Style = 'Normal'

$cell = $ws.Range("E14")
$cell.NumberFormat = '@'
$cell.Value = '0.90%'
$cell.Style = 'Normal'

$ws.Range("B15").Value = 'BitForexToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$cell = $ws.Range("D15")
$cell.NumberFormat = '@'
$cell.Value = '0.001281'
$cell.Style = 'Normal'

$cell = $ws.Range("E15")
$cell.NumberFormat = '@'
$cell.Value = '2.52%'
$cell.Style = 'Normal'

$ws.Range("B16").Value = 'TigerCash'

$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$cell = $ws.Range("D16")
$cell.NumberFormat = '@'
$cell.Value = '0.005884'
$cell.Style = 'Normal'

$cell = $ws.Range("E16")
$cell.NumberFormat = '@'
$cell.Value = '1.07%'
$cell.Style = 'Normal'

$ws.Range("B17").Value = 'LEO'

$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$cell = $ws.Range("D17")
$cell.NumberFormat = '@'
$cell.Value = '3.351'
$cell.Style = 'Normal'

$cell = $ws.Range("E17")
$cell.NumberFormat = '@'
$cell.Value = '-0.16%'
$cell.Style = 'Normal'

$cell = $ws.Range("D19")
$cell.NumberFormat = '@'
$cell.Value = '7.538'
$cell.Style = 'Normal'

$cell = $ws.Range("E19")
$cell.NumberFormat = '@'
$cell.Value = '12.49%'
$cell.Style = 'Normal'

$cell = $ws.Range("E20")
$cell.NumberFormat = '@'
$cell.Value = '-1.64%'
$cell.Style = 'Normal'

$cell = $ws.Range("D21")
$cell.NumberFormat = '@'
$cell.Value = '0.2719'
$cell.Style = 'Normal'

$cell = $ws.Range("E21")
$cell.NumberFormat = '@'
$cell.Value = '-0.58%'
$cell.Style = 'Normal'

$cell = $ws.Range("D22")
$cell.NumberFormat = '@'
$cell.Value = '0.04029'
$cell.Style = 'Normal'

$cell = $ws.Range("E22")
$cell.NumberFormat = '@'
$cell.Value = '-0.86%'
$cell.Style = 'Normal'

$cell = $ws.Range("E23")
$cell.NumberFormat = '@'
$cell.Value = '-0.51%'
$cell.Style = 'Normal'

$cell = $ws.Range("D24")
$cell.NumberFormat = '@'
$cell.Value = '0.004241'
$cell.Style = 'Normal'

$cell = $ws.Range("E24")
$cell.NumberFormat = '@'
$cell.Value = '4.30%'
$cell.Style = 'Normal'

$cell = $ws.Range("D25")
$cell.NumberFormat = '@'
$cell.Value = '0.0001270'
$cell.Style = 'Normal'

$cell = $ws.Range("E25")
$cell.NumberFormat = '@'
$cell.Value = '0.04%'
$cell.Style = 'Normal'

$cell = $ws.Range("D38")
$cell.NumberFormat = '@'
$cell.Value = '0.02551'
$cell.Style = 'Normal'

$cell = $ws.Range("E38")
$cell.NumberFormat = '@'
$cell.Value = '4.87%'
$cell.Style = 'Normal'

$cell = $ws.Range("E39")
$cell.NumberFormat = '@'
$cell.Value = '0.89%'
$cell.Style = 'Normal'

$cell = $ws.Range("D40")
$cell.NumberFormat = '@'
$cell.Value = '0.007835'
$cell.Style = 'Normal'

$cell = $ws.Range("E40")
$cell.NumberFormat = '@'
$cell.Value = '0.73%'
$cell.Style = 'Normal'

$cell = $ws.Range("D41")
$cell.NumberFormat = '@'
$cell.Value = '0.1316'
$cell.Style = 'Normal'

$cell = $ws.Range("E41")
$cell.NumberFormat = '@'
$cell.Value = '0.87%'
$cell.Style = 'Normal'

$cell = $ws.Range("D42")
$cell.NumberFormat = '@'
$cell.Value = '0.006641'
$cell.Style = 'Normal'

$cell = $ws.Range("E42")
$cell.NumberFormat = '@'
$cell.Value = '-2.18%'
$cell.Style = 'Normal'

$cell = $ws.Range("D43")
$cell.NumberFormat = '@'
$cell.Value = '0.001940'
$cell.Style = 'Normal'

$cell = $ws.Range("E43")
$cell.NumberFormat = '@'
$cell.Value = '-0.43%'
$cell.Style = 'Normal'

$cell = $ws.Range("D44")
$cell.NumberFormat = '@'
$cell.Value = '0.008096'
$cell.Style = 'Normal'

$cell = $ws.Range("E44")
$cell.NumberFormat = '@'
$cell.Value = '6.79%'
$cell.Style = 'Normal'

$cell = $ws.Range("D45")
$cell.NumberFormat = '@'
$cell.Value = '0.3094'
$cell.Style = 'Normal'

$cell = $ws.Range("E45")
$cell.NumberFormat = '@'
$cell.Value = '0.42%'
$cell.Style = 'Normal'

$cell = $ws.Range("D46")
$cell.NumberFormat = '@'
$cell.Value = '0.00006777'
$cell.Style = 'Normal'

$cell = $ws.Range("E46")
$cell.NumberFormat = '@'
$cell.Value = '-1.65%'
$cell.Style = 'Normal'

$cell = $ws.Range("D47")
$cell.NumberFormat = '@'
$cell.Value = '0.00000000750'
$cell.Style = 'Normal'

$cell = $ws.Range("E47")
$cell.NumberFormat = '@'
$cell.Value = '0.04%'
$cell.Style = 'Normal'

$cell = $ws.Range("D48")
$cell.NumberFormat = '@'
$cell.Value = '0.2181'
$cell.Style = 'Normal'

$cell = $ws.Range("E48")
$cell.NumberFormat = '@'
$cell.Value = '191.83%'
$cell.Style = 'Normal'

$cell = $ws.Range("E49")
$cell.NumberFormat = '@'
$cell.Value = '3.47%'
$cell.Style = 'Normal'

$cell = $ws.Range("D50")
$cell.NumberFormat = '@'
$cell.Value = '0.00002101'
$cell.Style = 'Normal'

$cell = $ws.Range("E50")
$cell.NumberFormat = '@'
$cell.Value = '0.04%'
$cell.Style = 'Normal'

$cell = $ws.Range("D51")
$cell.NumberFormat = '@'
$cell.Value = '0.0002001'
$cell.Style = 'Normal'

$cell = $ws.Range("E51")
$cell.NumberFormat = '@'
$cell.Value = '0.04%'
$cell.Style = 'Normal'
